$p = $ppt.ActivePresentation
$m = $p.SlideMaster
Write-Host "Master.Name before: " $m.Name
$m.Name = "Office Theme"
Write-Host "Master.Name after: " $m.Name
